$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "113.1 (#18)"
$ws.Range("E2").Value = "108.9 (#6)"
$ws.Range("C3").Value = "+4.2 (#8)"
$ws.Range("E3").Value = "-4.2 (#8)"
$ws.Range("C4").Value = "13.4 (#26)"
$ws.Range("E4").Value = "14.2 (#6)"
$ws.Range("C5").Value = "45.9 (#22)"
$ws.Range("C6").Value = "25.9 (#17)"
$ws.Range("E6").Value = "24.4 (#4)"
$ws.Range("C7").Value = "52.2 (#17)"
$ws.Range("E7").Value = "51.4 (#7)"
$ws.Range("C8").Value = "0.640 (#14)"
$ws.Range("E8").Value = "0.606 (#4)"
$ws.Range("C9").Value = "1.787 (#20)"
$ws.Range("E9").Value = "1.658 (#6)"
$ws.Range("C10").Value = "55.0% (#11)"
$ws.Range("E10").Value = "52.9% (#6)"
$ws.Range("C11").Value = "0.247 (#13)"
$ws.Range("E11").Value = "0.232 (#11)"
$ws.Range("C12").Value = "79.1% (#9)"
$ws.Range("E12").Value = "77.4% (#8)"
$ws.Range("C13").Value = "37.9% (#5)"
$ws.Range("E13").Value = "34.8% (#3)"
$ws.Range("C14").Value = "53.5% (#20)"
$ws.Range("E14").Value = "53.4% (#9)"
$ws.Range("C15").Value = "46.4% (#17)"
$ws.Range("C16").Value = "1.159 (#11)"
$ws.Range("C17").Value = "40.5 (#22)"
$ws.Range("E17").Value = "40.3 (#5)"
$ws.Range("C19").Value = "15.1 (#5)"
$ws.Range("C20").Value = "39.9 (#5)"
$ws.Range("E20").Value = "36.0 (#6)"
$ws.Range("C21").Value = "17.0 (#14)"
$ws.Range("E21").Value = "15.8 (#5)"
$ws.Range("C22").Value = "21.5 (#19)"
$ws.Range("E22").Value = "20.4 (#5)"
$ws.Range("C23").Value = "28.2 (#16)"
$ws.Range("E23").Value = "27.6 (#8)"
$ws.Range("C24").Value = "27.3 (#24)"
$ws.Range("E24").Value = "26.8 (#2)"
$ws.Range("C25").Value = "28.7 (#19)"
$ws.Range("C26").Value = "28.6 (#7)"
$ws.Range("C27").Value = "0.5 (#15)"
$ws.Range("E28").Value = "11.2 (#16)"
$ws.Range("C29").Value = "33.1 (#18)"
$ws.Range("E29").Value = "32.3 (#10)"
$ws.Range("C30").Value = "25.4% (#16)"
$ws.Range("E30").Value = "25.3% (#15)"
$ws.Range("C31").Value = "74.7% (#15)"
$ws.Range("E31").Value = "74.6% (#16)"
$ws.Range("E32").Value = "4.6 (#11)"
$ws.Range("C33").Value = "5.5% (#15)"
$ws.Range("E33").Value = "5.3% (#12)"
$ws.Range("C34").Value = "8.1 (#13)"
$ws.Range("E34").Value = "8.6 (#19)"
$ws.Range("C35").Value = "7.2% (#11)"
$ws.Range("E35").Value = "7.7% (#21)"
$ws.Range("C37").Value = "12.9% (#19)"
$ws.Range("C38").Value = "18.1 (#10)"
$ws.Range("E38").Value = "19.0 (#11)"
$ws.Range("E39").Value = "0.3 (#15)"
$ws.Range("C40").Value = "16.1% (#13)"
$ws.Range("E40").Value = "16.9% (#8)"
